# "variable prices and extra level"
# Applies to Perseverance/data/level1.xlsx:
#  - removes the unused trailing columns CW:DP (20 cols), shrinking the
#    used range from A1:DP13 to A1:CV13
#  - flips a handful of individual 0/1/2 cell values (price/level tweaks)
#  - updates the view (topLeftCell / selection)
#  - rebuilds the two "equals 1" / "equals 2" conditional-format rules so
#    four extra dxf records exist in styles.xml (count 2 -> 6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the trailing columns CW:DP -- nothing lives to their right, so
#    deleting them shrinks the used range to A1:CV13 exactly like the diff.
# ---------------------------------------------------------------------
$ws.Range("CW1:DP13").EntireColumn.Delete()

# ---------------------------------------------------------------------
# 2. Per-cell value edits (variable prices)
# ---------------------------------------------------------------------
$ws.Range("BY1").Value = 0
$ws.Range("CT1").Value = 1
$ws.Range("CU1").Value = 1
$ws.Range("BY2").Value = 0
$ws.Range("CT2").Value = 1
$ws.Range("CU2").Value = 1
$ws.Range("CT3").Value = 1
$ws.Range("CU3").Value = 1
$ws.Range("CT4").Value = 1
$ws.Range("CU4").Value = 1
$ws.Range("CT5").Value = 1
$ws.Range("CU5").Value = 1
$ws.Range("BY6").Value = 0
$ws.Range("CU6").Value = 1
$ws.Range("CH7").Value = 2
$ws.Range("CM7").Value = 2
$ws.Range("CS7").Value = 2
$ws.Range("CT7").Value = 2
$ws.Range("CC8").Value = 2
$ws.Range("CS8").Value = 2
$ws.Range("CT8").Value = 2
$ws.Range("CI9").Value = 1
$ws.Range("CJ9").Value = 1
$ws.Range("CK9").Value = 1
$ws.Range("CL9").Value = 1
$ws.Range("CN9").Value = 1
$ws.Range("CO9").Value = 1
$ws.Range("CP9").Value = 1
$ws.Range("CQ9").Value = 1
$ws.Range("CS9").Value = 2
$ws.Range("CT9").Value = 2
$ws.Range("CD10").Value = 1
$ws.Range("CE10").Value = 1
$ws.Range("CF10").Value = 1
$ws.Range("CG10").Value = 1
$ws.Range("CL10").Value = 0
$ws.Range("CM10").Value = 0
$ws.Range("CN10").Value = 0
$ws.Range("CO10").Value = 0
$ws.Range("CQ10").Value = 0
$ws.Range("CR10").Value = 0
$ws.Range("CS10").Value = 0
$ws.Range("CT10").Value = 0
$ws.Range("CU10").Value = 1
$ws.Range("BY11").Value = 1
$ws.Range("CC11").Value = 0
$ws.Range("CF11").Value = 0
$ws.Range("CG11").Value = 0
$ws.Range("CH11").Value = 0
$ws.Range("CI11").Value = 0
$ws.Range("CJ11").Value = 0
$ws.Range("CT11").Value = 1
$ws.Range("CU11").Value = 1
$ws.Range("BX12").Value = 0
$ws.Range("CT12").Value = 1
$ws.Range("CU12").Value = 1
$ws.Range("BX13").Value = 0
$ws.Range("CT13").Value = 1
$ws.Range("CU13").Value = 1

# ---------------------------------------------------------------------
# 3. View state: scrolled one column further right, new selection
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 36
$ws.Range("CK42").Select()

# ---------------------------------------------------------------------
# 4. Rebuild conditional formatting so new dxf records are generated
# ---------------------------------------------------------------------
$old = $ws.Range("A1:FD13")
$old.FormatConditions.Delete()

$full = $ws.Range("A1:CV13")
$cf1 = $full.FormatConditions.Add(1, 3, "=2")
$cf1.Interior.Color = 10284031
$cf1.Font.Color = 26012

$cf2 = $full.FormatConditions.Add(1, 3, "=1")
$cf2.Interior.Color = 13551615

Write-Host "done"
